$d = $word.ActiveDocument

$old = "年キャンペーン期間 対象：ペルセウス座 2022: 1月16日〜25日、11月7日〜16日、12月6日〜15日"
$new = " ：2022年キャンペーン期間 (対象：ペルセウス座)：、1月16日〜25日、11月7日〜16日、12月6日〜15日"

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
